# Add a "Kaggle link: <url>" paragraph to the "Tools used" slide (slide 14),
# just before the trailing empty paragraph in the content placeholder.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(14)
$shp = $s.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange

$kaggleUrl = "https://www.kaggle.com/c/house-prices-advanced-regression-techniques/"

# The content placeholder currently has 5 paragraphs: Python / Jupyter Notebook /
# Libraries.../ Visualization.../ (trailing empty bullet-less paragraph).
# Insert a new paragraph right after paragraph 4 ("Visualization: Tableau, PowerPoint")
# containing "Kaggle link: " followed by the URL run, before the trailing paragraph.
$visualizationPara = $tr.Paragraphs(4)
$null = $visualizationPara.InsertAfter([char]13 + "Kaggle link: ")

$newPara = $tr.Paragraphs(5)
$null = $newPara.InsertAfter($kaggleUrl)

# Re-fetch the new paragraph and compute the precise character span of the URL
# text within it, so formatting/hyperlink only ever touches that run (and not
# the "Kaggle link: " prefix run).
$p5 = $tr.Paragraphs(5)
$prefixLen = "Kaggle link: ".Length
$urlLen = $kaggleUrl.Length
$urlStart = $p5.Start + $prefixLen

# Apply the hyperlink first (on a freshly scoped range covering only the URL
# text), then shrink its font size - this ordering keeps the paragraph-end
# run properties ("endParaRPr") at the paragraph's base size (14pt) instead
# of inheriting the smaller hyperlink run size.
$urlForLink = $tr.Characters($urlStart, $urlLen)
$ppMouseClick = 1
$action = $urlForLink.ActionSettings($ppMouseClick)
$action.Hyperlink.Address = $kaggleUrl

$urlForSize = $tr.Characters($urlStart, $urlLen)
$urlForSize.Font.Size = 10.5
